# This workbook is backed by a SharePoint ("owssvr") list query. The list
# was refreshed upstream: the "Adobe eLearning" item was removed from the
# list, and its distinguishing taxonomy tags were folded into the
# "Adobe - Captivate" item. After the refresh, the query/table was renamed
# (Excel appends a "(4)" / "__4" suffix because it is the 4th refresh of
# this query saved under a new name) and the data range shrank by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fold the now-removed "Adobe eLearning" row's unique tag fragments into
# the "Adobe - Captivate" row (row 3) before the row disappears.
$ws.Range("B3").Value = "_Starting Over;#32;#Simulation;#40;#Learning;#37"
$ws.Range("C3").Value = "_Starting Over;#22;#Video;#45;#Online Module;#37;#Training Manual;#39"

# Remove the "Adobe eLearning" row entirely; every following row shifts up.
$ws.Rows(6).Delete()

# Rename the table/list object and its query table to match the refreshed
# query name ("owssvr (4)" / "Table_owssvr__4").
$lo = $ws.ListObjects.Item(1)
$lo.Name = "Table_owssvr__4"

try {
    $qt = $lo.QueryTable
    $qt.Name = "owssvr (4)"
} catch {
}

try {
    $conn = $wb.Connections.Item(1)
    $conn.Name = "owssvr (4)"
} catch {
}

# Update the hidden defined name used by the query table to reflect the
# new name and the now-smaller data range.
foreach ($n in $wb.Names) {
    $n.Name = "owssvr__4"
    $n.RefersTo = "=owssvr!`$A`$1:`$F`$32"
}

# Reposition/resize the saved window to match the refreshed workbook view.
try {
    $w = $wb.Windows.Item(1)
    $w.Left = 120
    $w.Top = 150
    $w.Width = 24915
    $w.Height = 12075
} catch {
}
